# Append four new "xiaomi" price observations to the price-history sheet
# (mirrors the existing iphone/samsung rows: product, price($), date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("xiaomi", 802, 44562),
    @("xiaomi", 800, 44563),
    @("xiaomi", 750, 44564),
    @("xiaomi", 740, 44565)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]

    $dateCell = $ws.Cells.Item($r, 3)
    $dateCell.Value = $entry[2]
    $dateCell.NumberFormat = "mm/dd/yy"
}

# Match the workbook's recorded selection after the edit.
$ws.Range("C23").Select()
